$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay text-typed (avoid Excel auto-numeric coercion)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "22.413.37"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.563.23"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "1.000"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").Value = "285.50"
$ws.Range("D7").Value = "0.3636"
$ws.Range("E7").Value = "  -2.74%  "
$ws.Range("D8").Value = "48.49"
$ws.Range("E8").Value = "  -2.85%  "
$ws.Range("D9").Value = "0.3338"
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("D10").Value = "1.128"
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("D11").Value = "0.07414"
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("D13").Value = "20.81"
$ws.Range("E13").Value = "  -2.63%  "
$ws.Range("D14").Value = "5.930"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").Value = "6.892"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").Value = "1.564.57"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").Value = "0.00001106"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").Value = "88.24"
$ws.Range("E18").Value = "  -3.11%  "
$ws.Range("D19").Value = "0.06681"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "6.364"
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("D22").Value = "16.12"
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("D23").Value = "11.98"
$ws.Range("E23").Value = "  -1.24%  "
$ws.Range("D24").Value = "22.400.03"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").Value = "2.412"
$ws.Range("E25").Value = "  +3.52%  "
$ws.Range("D26").Value = "2.568"
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("D27").Value = "149.83"
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("E28").Value = "  -3.62%  "
$ws.Range("D29").Value = "4.995"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "123.06"
$ws.Range("E30").Value = "  -2.22%  "
$ws.Range("D31").Value = "1.739.36"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("D33").Value = "6.147"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("D34").Value = "1.993"
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("D35").Value = "9.807"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "0.08249"
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("D37").Value = "0.02399"
$ws.Range("E37").Value = "  -2.69%  "
$ws.Range("D38").Value = "1.307"
$ws.Range("E38").Value = "  -5.79%  "
$ws.Range("D39").Value = "0.06386"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("D40").Value = "0.2209"
$ws.Range("E40").Value = "  -3.45%  "
$ws.Range("D41").Value = "5.336"
$ws.Range("E41").Value = "  -2.23%  "
$ws.Range("D42").Value = "11.18"
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").Value = "0.6088"
$ws.Range("E43").Value = "  -2.59%  "
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("D46").Value = "3.758"
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("D47").Value = "0.5762"
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").Value = "2.012"
$ws.Range("E48").Value = "  -3.37%  "
$ws.Range("D49").Value = "124.89"
$ws.Range("E49").Value = "  -3.40%  "
$ws.Range("D50").Value = "1.215"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").Value = "0.07213"
$ws.Range("E51").Value = "  -1.55%  "
